$wb = $excel.ActiveWorkbook

# Sheet 1: 展览 (worksheet index 1)
$ws1 = $wb.Worksheets.Item(1)
$ws1.Cells.Item(2, 6).Value = 143
$ws1.Cells.Item(3, 6).Value = 412
$ws1.Cells.Item(4, 6).Value = 201
$ws1.Cells.Item(5, 6).Value = 43
$ws1.Cells.Item(6, 6).Value = 1263
$ws1.Cells.Item(7, 6).Value = 465
$ws1.Cells.Item(9, 6).Value = 218
$ws1.Cells.Item(11, 6).Value = 190
$ws1.Cells.Item(12, 6).Value = 1066
$ws1.Cells.Item(15, 6).Value = 212
$ws1.Cells.Item(16, 6).Value = 1553
$ws1.Cells.Item(17, 6).Value = 573
$ws1.Cells.Item(18, 6).Value = 244
$ws1.Cells.Item(19, 6).Value = 367
$ws1.Cells.Item(21, 6).Value = 872
$ws1.Cells.Item(22, 6).Value = 1173
$ws1.Cells.Item(24, 6).Value = 1911
$ws1.Cells.Item(25, 6).Value = 2704
$ws1.Cells.Item(26, 6).Value = 1492
$ws1.Cells.Item(28, 6).Value = 66
$ws1.Cells.Item(29, 6).Value = 486
$ws1.Cells.Item(30, 6).Value = 821
$ws1.Cells.Item(31, 6).Value = 1391
$ws1.Cells.Item(32, 6).Value = 842
$ws1.Cells.Item(33, 6).Value = 1471
$ws1.Cells.Item(34, 6).Value = 172
$ws1.Cells.Item(36, 6).Value = 801
$ws1.Cells.Item(37, 6).Value = 680
$ws1.Cells.Item(38, 6).Value = 702
$ws1.Cells.Item(39, 6).Value = 897
$ws1.Cells.Item(40, 6).Value = 378
$ws1.Cells.Item(41, 6).Value = 270

# Sheet 2: 演出 (worksheet index 2)
$ws2 = $wb.Worksheets.Item(2)
$ws2.Cells.Item(3, 6).Value = 3
$ws2.Cells.Item(5, 6).Value = 8
$ws2.Cells.Item(15, 6).Value = 692
$ws2.Cells.Item(19, 6).Value = 14

# Sheet 4: 全部类型 (worksheet index 4)
$ws4 = $wb.Worksheets.Item(4)
$ws4.Cells.Item(3, 6).Value = 143
$ws4.Cells.Item(4, 6).Value = 412
$ws4.Cells.Item(5, 6).Value = 201
$ws4.Cells.Item(6, 6).Value = 43
$ws4.Cells.Item(7, 6).Value = 8
$ws4.Cells.Item(9, 6).Value = 1263
$ws4.Cells.Item(10, 6).Value = 465
$ws4.Cells.Item(12, 6).Value = 218
$ws4.Cells.Item(14, 6).Value = 190
$ws4.Cells.Item(15, 6).Value = 1066
$ws4.Cells.Item(18, 6).Value = 212
$ws4.Cells.Item(19, 6).Value = 1554
$ws4.Cells.Item(20, 6).Value = 573
$ws4.Cells.Item(21, 6).Value = 244
$ws4.Cells.Item(22, 6).Value = 367
$ws4.Cells.Item(25, 6).Value = 1173
$ws4.Cells.Item(26, 6).Value = 2704
$ws4.Cells.Item(28, 6).Value = 1492
$ws4.Cells.Item(31, 6).Value = 66
$ws4.Cells.Item(34, 6).Value = 487
$ws4.Cells.Item(35, 6).Value = 821
$ws4.Cells.Item(36, 6).Value = 1391
$ws4.Cells.Item(39, 6).Value = 842
$ws4.Cells.Item(40, 6).Value = 1472
$ws4.Cells.Item(41, 6).Value = 801
$ws4.Cells.Item(42, 6).Value = 680
$ws4.Cells.Item(43, 6).Value = 702
$ws4.Cells.Item(44, 6).Value = 897
$ws4.Cells.Item(45, 6).Value = 378
$ws4.Cells.Item(48, 6).Value = 270
